# Actualización automática 2025-10-21 13:30:08
# Inserts a new advisor/client row ("ALMENDARIZ MOLINA HENRRY MAURICIO") at
# row 2 (alphabetically before "ALTAMIRANO ...") on both the
# "VENTAS POR GRUPO" and "VENTA MENSUAL" sheets, shifting all subsequent
# rows (and the trailing totals row) down by one.

$wb = $excel.ActiveWorkbook

$currencyFormat = '"$"#,##0.00'
$newClient = "ALMENDARIZ MOLINA HENRRY MAURICIO"
$asesor = "ALMEIDA CUATIN JHONATHANN CARLOS"

# ---------------------------------------------------------------------
# Sheet 1: "VENTAS POR GRUPO"  (columns A:R, data rows 2-35 -> 3-36,
# totals row 36 -> 37)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

$ws1.Rows.Item(2).Insert()

$ws1.Range("A2:R2").ClearFormats()
$ws1.Range("C2:R2").NumberFormat = $currencyFormat

$ws1.Range("A2").Value = $asesor
$ws1.Range("B2").Value = $newClient
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(2, $c).Value = 0
}

# Update the "X de 34" -> "X de 35" counters on the totals row, which is
# now row 37 after the insert shifted it down from row 36.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(37, $c)
    $label = $cell.Text
    $cell.Value = $label -replace "de 34", "de 35"
}

# ---------------------------------------------------------------------
# Sheet 2: "VENTA MENSUAL"  (columns A:G, data rows 2-35 -> 3-36,
# totals row 36 -> 37)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Rows.Item(2).Insert()

$ws2.Range("A2:G2").ClearFormats()
$ws2.Range("C2:G2").NumberFormat = $currencyFormat

$ws2.Range("A2").Value = $asesor
$ws2.Range("B2").Value = $newClient
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(2, $c).Value = 0
}
